$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# Update the time_taken (column F) timestamps on the "data" sheet
$newTimes = @(
    "2021-10-05 14:34:56.318911",
    "2021-10-05 14:34:56.318919",
    "2021-10-05 14:34:56.318922",
    "2021-10-05 14:34:56.318925",
    "2021-10-05 14:34:56.318928",
    "2021-10-05 14:34:56.318931",
    "2021-10-05 14:34:56.318934",
    "2021-10-05 14:34:56.318936",
    "2021-10-05 14:34:56.318939",
    "2021-10-05 14:34:56.318942",
    "2021-10-05 14:34:56.318944",
    "2021-10-05 14:34:56.318947",
    "2021-10-05 14:34:56.318950",
    "2021-10-05 14:34:56.318952",
    "2021-10-05 14:34:56.318955",
    "2021-10-05 14:34:56.318957",
    "2021-10-05 14:34:56.318960",
    "2021-10-05 14:34:56.318963",
    "2021-10-05 14:34:56.318966",
    "2021-10-05 14:34:56.318969",
    "2021-10-05 14:34:56.318971",
    "2021-10-05 14:34:56.318974",
    "2021-10-05 14:34:56.318976",
    "2021-10-05 14:34:56.318979",
    "2021-10-05 14:34:56.318981",
    "2021-10-05 14:34:56.318984",
    "2021-10-05 14:34:56.318987",
    "2021-10-05 14:34:56.318989",
    "2021-10-05 14:34:56.318992",
    "2021-10-05 14:34:56.318994",
    "2021-10-05 14:34:56.318997",
    "2021-10-05 14:34:56.318999",
    "2021-10-05 14:34:56.319002",
    "2021-10-05 14:34:56.319005",
    "2021-10-05 14:34:56.319007",
    "2021-10-05 14:34:56.319010",
    "2021-10-05 14:34:56.319012",
    "2021-10-05 14:34:56.319015",
    "2021-10-05 14:34:56.319017",
    "2021-10-05 14:34:56.319020",
    "2021-10-05 14:34:56.319023",
    "2021-10-05 14:34:56.319026",
    "2021-10-05 14:34:56.319028",
    "2021-10-05 14:34:56.319031",
    "2021-10-05 14:34:56.319033",
    "2021-10-05 14:34:56.319036",
    "2021-10-05 14:34:56.319038",
    "2021-10-05 14:34:56.319041",
    "2021-10-05 14:34:56.319043",
    "2021-10-05 14:34:56.319046",
    "2021-10-05 14:34:56.319048",
    "2021-10-05 14:34:56.319051",
    "2021-10-05 14:34:56.319054",
    "2021-10-05 14:34:56.319057",
    "2021-10-05 14:34:56.319059",
    "2021-10-05 14:34:56.319062",
    "2021-10-05 14:34:56.319064",
    "2021-10-05 14:34:56.319067",
    "2021-10-05 14:34:56.319069",
    "2021-10-05 14:34:56.319072",
    "2021-10-05 14:34:56.319074",
    "2021-10-05 14:34:56.319077",
    "2021-10-05 14:34:56.319079",
    "2021-10-05 14:34:56.319082",
    "2021-10-05 14:34:56.319086",
    "2021-10-05 14:34:56.319088"
)
for ($i = 0; $i -lt $newTimes.Count; $i++) {
    $row = $i + 2
    $dataSheet.Cells.Item($row, 6).Value = $newTimes[$i]
}

# Add a new "metadata" worksheet positioned after "data"
$metaSheet = $wb.Worksheets.Add($null, $dataSheet)
$metaSheet.Name = "metadata"

# Header row (row 1), matching the bold/bordered header style used on "data"
$headerStyle = $dataSheet.Range("B1").Style
$metaSheet.Range("B1").Value = "data_name"
$metaSheet.Range("C1").Value = "data_id"
$metaSheet.Range("D1").Value = "data_version"
$metaSheet.Range("E1").Value = "data_version_created"
$metaSheet.Range("F1").Value = "panel_query_time"
$metaSheet.Range("G1").Value = "panel_get_request"
$metaSheet.Range("B1:G1").Style = $headerStyle

# Data row (row 2)
$indexStyle = $dataSheet.Range("A2").Style
$metaSheet.Range("A2").Value = 0
$metaSheet.Range("A2").Style = $indexStyle
$metaSheet.Range("B2").Value = "Motor Neurone Disease"
$metaSheet.Range("C2").Value = 25
$metaSheet.Range("D2").Value = "'0.131"
$metaSheet.Range("E2").Value = "2021-08-29T08:12:13.011006Z"
$metaSheet.Range("F2").Value = "2021-10-05 14:34:56.315230"
$metaSheet.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/25/?format=json"

